$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 already carries the "Normal"-ish style (cellXfs index 1) that the rest
# of the header row uses; copy its format and stamp it onto the cells the
# diff touches so they end up sharing that same style index.
$styleSource = $ws.Range("A1")
$styleSource.Copy()

# Row 1: M1 / N1 keep their existing text (es / he_IL), just gain the style.
foreach ($ref in @("M1", "N1")) {
    $ws.Range($ref).PasteSpecial(-4122)  # xlPasteFormats
}

# Rows 2-4: columns D, F, H, I, K, L, M, N all become the literal "test"
# string, styled the same way. Column J (the date) is left untouched.
$cols = @("D", "F", "H", "I", "K", "L", "M", "N")
foreach ($row in 2..4) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $cell.Value = "test"
        $cell.PasteSpecial(-4122)  # xlPasteFormats
    }
}

$excel.CutCopyMode = 0
